$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row of data (2025-10-23) below the existing last row (66).
# Use a leading apostrophe so the date-like text is stored literally (as a
# string) instead of being auto-converted into a date serial number, then
# reset the style back to Normal so no extra (quote-prefixed) number format
# is left behind on the cell - matching how the other date-text rows above
# it are stored (plain inline/shared string, default style).
$ws.Range("A67").Value = "'10/23/2025"
$ws.Range("A67").Style = "Normal"

$ws.Range("B67").Value = 10155.31
